# Auto-generated edit script applying the diff to Asura_Profits workbook
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1846.3846
$ws.Range("I132").Value = 1581.3438
$ws.Range("J132").Value = 3058
$ws.Range("K132").Value = 4744.0314
$ws.Range("L132").Value = 9174
$ws.Range("M132").Value = -2214.0314
$ws.Range("N132").Value = -14234

$ws.Range("H138").Value = 3188.8245
$ws.Range("I138").Value = 1412.5172
$ws.Range("K138").Value = 4237.5516
$ws.Range("M138").Value = 902.4484000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 0
$ws.Range("I6").Value = 0
$ws.Range("K6").Value = 0
$ws.Range("M6").Value = $null

$ws.Range("H61").Value = 2528.1143
$ws.Range("I61").Value = 2413.3103
$ws.Range("J61").Value = 3083
$ws.Range("K61").Value = 2413.3103
$ws.Range("L61").Value = 3083
$ws.Range("M61").Value = -2201.3103
$ws.Range("N61").Value = -3507

$ws.Range("H135").Value = 39851.145
$ws.Range("J135").Value = 39851.145
$ws.Range("L135").Value = 39851.145
$ws.Range("N135").Value = -49991.145

$ws.Range("H136").Value = 2528.1143
$ws.Range("I136").Value = 2413.3103
$ws.Range("J136").Value = 3083
$ws.Range("K136").Value = 7239.9309
$ws.Range("L136").Value = 9249
$ws.Range("M136").Value = -4689.9309
$ws.Range("N136").Value = -14349

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H7").Value = 16512.143
$ws.Range("I7").Value = 6920.75
$ws.Range("K7").Value = 6920.75
$ws.Range("M7").Value = -6807.75

$ws.Range("H134").Value = 2236.8696
$ws.Range("I134").Value = 2207.4
$ws.Range("J134").Value = 2433.3333
$ws.Range("K134").Value = 6622.200000000001
$ws.Range("L134").Value = 7299.999899999999
$ws.Range("M134").Value = -4087.200000000001
$ws.Range("N134").Value = -12369.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H3").Value = 17700
$ws.Range("J3").Value = 14625
$ws.Range("L3").Value = 14625
$ws.Range("N3").Value = -14851

$ws.Range("H7").Value = 186.4
$ws.Range("I7").Value = 65
$ws.Range("J7").Value = 267.33334
$ws.Range("K7").Value = 65
$ws.Range("L7").Value = 267.33334
$ws.Range("M7").Value = 48
$ws.Range("N7").Value = -493.33334

$ws.Range("H17").Value = 20000
$ws.Range("J17").Value = 20000
$ws.Range("L17").Value = 20000
$ws.Range("N17").Value = -20348

$ws.Range("H22").Value = 268.90475
$ws.Range("I22").Value = 248.70589
$ws.Range("J22").Value = 354.75
$ws.Range("K22").Value = 248.70589
$ws.Range("L22").Value = 354.75
$ws.Range("M22").Value = 101.29411
$ws.Range("N22").Value = -1054.75

$ws.Range("H31").Value = 1772.8628
$ws.Range("I31").Value = 1395.1
$ws.Range("J31").Value = 2312.524
$ws.Range("K31").Value = 1395.1
$ws.Range("L31").Value = 2312.524
$ws.Range("M31").Value = -1100.1
$ws.Range("N31").Value = -2902.524

$ws.Range("H34").Value = 1772.8628
$ws.Range("I34").Value = 1395.1
$ws.Range("J34").Value = 2312.524
$ws.Range("K34").Value = 1395.1
$ws.Range("L34").Value = 2312.524
$ws.Range("M34").Value = -1193.1
$ws.Range("N34").Value = -2716.524

$ws.Range("H107").Value = 412.48486
$ws.Range("I107").Value = 383.875
$ws.Range("J107").Value = 488.77777
$ws.Range("K107").Value = 383.875
$ws.Range("L107").Value = 488.77777
$ws.Range("M107").Value = 1536.125
$ws.Range("N107").Value = -4328.77777

$ws.Range("H134").Value = 1824.6316
$ws.Range("I134").Value = 1457.0667
$ws.Range("K134").Value = 4371.2001
$ws.Range("M134").Value = -1836.2001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 86758.086
$ws.Range("J4").Value = 6766.6665
$ws.Range("L4").Value = 20299.9995
$ws.Range("N4").Value = -20523.9995

$ws.Range("H12").Value = 172.5
$ws.Range("I12").Value = 95
$ws.Range("J12").Value = 250
$ws.Range("K12").Value = 285
$ws.Range("L12").Value = 750
$ws.Range("M12").Value = -112
$ws.Range("N12").Value = -1096

$ws.Range("H17").Value = 3400
$ws.Range("J17").Value = 5500
$ws.Range("L17").Value = 16500
$ws.Range("N17").Value = -16838

$ws.Range("H68").Value = 1115.2106
$ws.Range("I68").Value = 0
$ws.Range("J68").Value = 1115.2106
$ws.Range("K68").Value = 0
$ws.Range("L68").Value = 3345.6318
$ws.Range("M68").Value = $null
$ws.Range("N68").Value = -4967.6318

$ws.Range("H71").Value = 1115.2106
$ws.Range("I71").Value = 0
$ws.Range("J71").Value = 1115.2106
$ws.Range("K71").Value = 0
$ws.Range("L71").Value = 10036.8954
$ws.Range("M71").Value = $null
$ws.Range("N71").Value = -18148.8954

$ws.Range("H112").Value = 4169.909
$ws.Range("I112").Value = 2658.1667
$ws.Range("J112").Value = 5984
$ws.Range("K112").Value = 7974.500100000001
$ws.Range("L112").Value = 17952
$ws.Range("M112").Value = -6866.500100000001
$ws.Range("N112").Value = -20168

$ws.Range("H113").Value = 294838.88
$ws.Range("I113").Value = 833877.75
$ws.Range("J113").Value = 817.6818
$ws.Range("K113").Value = 2501633.25
$ws.Range("L113").Value = 2453.0454
$ws.Range("M113").Value = -2499463.25
$ws.Range("N113").Value = -6793.0454

$ws.Range("H122").Value = 1371.3572
$ws.Range("I122").Value = 600
$ws.Range("J122").Value = 1430.6923
$ws.Range("K122").Value = 5400
$ws.Range("L122").Value = 12876.2307
$ws.Range("M122").Value = -2950
$ws.Range("N122").Value = -17776.2307

$ws.Range("H131").Value = 3274.7234
$ws.Range("I131").Value = 665.5714
$ws.Range("J131").Value = 3731.325
$ws.Range("K131").Value = 1996.7142
$ws.Range("L131").Value = 11193.975
$ws.Range("M131").Value = 3043.2858
$ws.Range("N131").Value = -21273.975

$ws.Range("H133").Value = 5404.4443
$ws.Range("I133").Value = 2560
$ws.Range("J133").Value = 6826.6665
$ws.Range("K133").Value = 7680
$ws.Range("L133").Value = 20479.9995
$ws.Range("M133").Value = -2620
$ws.Range("N133").Value = -30599.9995

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 28575.715

$ws.Range("H70").Value = 7880.24
$ws.Range("I70").Value = 7958.5835
$ws.Range("J70").Value = 6000
$ws.Range("K70").Value = 7958.5835
$ws.Range("L70").Value = 6000
$ws.Range("M70").Value = -7688.5835
$ws.Range("N70").Value = -6540

$ws.Range("H73").Value = 7880.24
$ws.Range("I73").Value = 7958.5835
$ws.Range("J73").Value = 6000
$ws.Range("K73").Value = 7958.5835
$ws.Range("L73").Value = 6000
$ws.Range("M73").Value = -7022.5835
$ws.Range("N73").Value = -7872

$ws.Range("H80").Value = 3185.6667
$ws.Range("I80").Value = 2991.3635
$ws.Range("J80").Value = 3720
$ws.Range("K80").Value = 2991.3635
$ws.Range("L80").Value = 3720
$ws.Range("M80").Value = -1993.3635
$ws.Range("N80").Value = -5716

$ws.Range("H83").Value = 3185.6667
$ws.Range("I83").Value = 2991.3635
$ws.Range("J83").Value = 3720
$ws.Range("K83").Value = 14956.8175
$ws.Range("L83").Value = 18600
$ws.Range("M83").Value = -9964.817499999999
$ws.Range("N83").Value = -28584

$ws.Range("H102").Value = 2936.7036
$ws.Range("I102").Value = 2617.4119
$ws.Range("J102").Value = 3479.5
$ws.Range("K102").Value = 2617.4119
$ws.Range("L102").Value = 3479.5
$ws.Range("M102").Value = -995.4119000000001
$ws.Range("N102").Value = -6723.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4142.7188
$ws.Range("I132").Value = 4152.2
$ws.Range("J132").Value = 4000.5
$ws.Range("K132").Value = 12456.6
$ws.Range("L132").Value = 12001.5
$ws.Range("M132").Value = -9926.599999999999
$ws.Range("N132").Value = -17061.5

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H39").Value = 13157980
$ws.Range("I39").Value = 39393940
$ws.Range("J39").Value = 40000
$ws.Range("K39").Value = 39393940
$ws.Range("L39").Value = 40000
$ws.Range("M39").Value = -39393527
$ws.Range("N39").Value = -40826

$ws.Range("H42").Value = 43808.8
$ws.Range("J42").Value = 39666.668
$ws.Range("L42").Value = 39666.668
$ws.Range("N42").Value = -40422.668

$ws.Range("H43").Value = 0
$ws.Range("I43").Value = 0
$ws.Range("K43").Value = 0
$ws.Range("M43").Value = $null

$ws.Range("H68").Value = 29925
$ws.Range("J68").Value = 29925
$ws.Range("L68").Value = 29925
$ws.Range("N68").Value = -31547

$ws.Range("H71").Value = 29925
$ws.Range("J71").Value = 29925
$ws.Range("L71").Value = 89775
$ws.Range("N71").Value = -97887
